$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Story")

# --- Update existing rows: Design -> Implement / Deferred (E column) ---
$ws.Range("E13").Value = "Implement"
$ws.Range("E21").Value = "Implement"
$ws.Range("E22").Value = "Implement"
$ws.Range("E23").Value = "Implement"
$ws.Range("E38").Value = "Deferred"
$ws.Range("E39").Value = "Implement"
$ws.Range("E41").Value = "Implement"
$ws.Range("E42").Value = "Implement"
$ws.Range("E44").Value = "Implement"
$ws.Range("E45").Value = "Deferred"

# --- Append new user stories (rows 46-52) ---
$ws.Range("A46").Value = "S052"
$ws.Range("B46").Value = "Create new quote from Organization/Lead Detail View."
$ws.Range("C46").Value = "Logged In User"
$ws.Range("D46").Value = "Create a quote from the DetailView of Organization/Lead under the Quotes tab."
$ws.Range("E46").Value = "Design"
$ws.Range("F46").Value = "Could"
$ws.Range("J46").Value = "V 1.0"

$ws.Range("A47").Value = "S053"
$ws.Range("B47").Value = "Add a site address"
$ws.Range("C47").Value = "Logged In User"
$ws.Range("D47").Value = "While creating the Lead/Organization, divide the address box into 2 parts (left and right) to record the site address and the billing address."
$ws.Range("E47").Value = "Design"
$ws.Range("F47").Value = "Could"
$ws.Range("J47").Value = "V 1.0"

$ws.Range("A48").Value = "S053"
$ws.Range("B48").Value = "Before the Quote is converted into a Sales Order, capture the Negotiated Final Amount"
$ws.Range("C48").Value = "Logged In User"
$ws.Range("D48").Value = "Before the Quote is converted into a Sales Order, capture the Negotiated Final Amount"
$ws.Range("E48").Value = "Design"
$ws.Range("F48").Value = "Should"
$ws.Range("J48").Value = "V 1.0"

$ws.Range("A49").Value = "S054"
$ws.Range("B49").Value = "Add a new Contracter."
$ws.Range("C49").Value = "Logged In User"
$ws.Range("D49").Value = "Add a new Contracter."
$ws.Range("E49").Value = "Design"
$ws.Range("F49").Value = "Should"
$ws.Range("J49").Value = "V 1.0"

$ws.Range("A50").Value = "S055"
$ws.Range("B50").Value = "View the Line Items assignments."
$ws.Range("C50").Value = "Logged In User"
$ws.Range("D50").Value = "View the line items assignments to external contracters, if any, in the tab shown in the Quote DetailView."
$ws.Range("E50").Value = "Design"
$ws.Range("F50").Value = "Should"
$ws.Range("J50").Value = "V 1.0"

$ws.Range("A51").Value = "S056"
$ws.Range("B51").Value = "View the Order Status History"
$ws.Range("C51").Value = "Logged In User"
$ws.Range("D51").Value = "View the status changes for the Order."
$ws.Range("E51").Value = "Design"
$ws.Range("F51").Value = "Should"
$ws.Range("J51").Value = "V 1.0"

$ws.Range("A52").Value = "S057"
$ws.Range("B52").Value = "View the maintenance work orders in a different tab."
$ws.Range("C52").Value = "Logged In User"
$ws.Range("D52").Value = "View the Maintenance Work Orders in a different tab other that order."
$ws.Range("E52").Value = "Design"
$ws.Range("F52").Value = "Must"
$ws.Range("J52").Value = "V 1.0"

# --- Row heights for wrapped multi-line text (matches author's saved heights) ---
$ws.Rows.Item(46).RowHeight = 25.5
$ws.Rows.Item(47).RowHeight = 51
$ws.Rows.Item(48).RowHeight = 38.25
$ws.Rows.Item(50).RowHeight = 38.25
$ws.Rows.Item(52).RowHeight = 25.5

# --- View state: scroll + selection moved to reflect the newly entered rows ---
$ws.Application.ActiveWindow.ScrollRow = 39
$ws.Range("A53:D53").Select()
